$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Formatting fix-up -----------------------------------------------------
# The header cell A1 ("Greeeting") and the "bye" row (A3) currently carry the
# "plain" cell style, while every other populated cell in A1:D4 carries the
# Arial/theme-1 style. Re-balance the styling so that A3 (which already used
# the plain style) picks up A1's Arial/theme-1 look-alike group, while the
# rest of the data cells pick up the plain style that A3 used to have. Using
# Copy/PasteSpecial(formats) (rather than touching Font.* directly) makes
# Excel re-use the existing style records instead of minting new ones.

# 1) Push A3's current ("plain") format onto the rest of the data cells.
$ws.Range("A3").Copy() | Out-Null
$ws.Range("B1:C1").PasteSpecial(-4122) | Out-Null
$ws.Range("A2:D2").PasteSpecial(-4122) | Out-Null
$ws.Range("B3:C3").PasteSpecial(-4122) | Out-Null
$ws.Range("A4:C4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# 2) Push A1's original (Arial/theme-1) format onto A3.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("A3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Text fix-up -------------------------------------------------------------
# Correct the "Greeeting" typo in the header cell.
$ws.Range("A1").Value = "Greeting"
